$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.496.34"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "2.358.25"
$ws.Range("E3").Value = "  -5.46%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'513.00"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").Value = "'127.46"
$ws.Range("E6").Value = "  -5.81%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "2.373.22"
$ws.Range("E9").Value = "  -5.83%  "
$ws.Range("D10").Value = "'0.0957"
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "'4.81"
$ws.Range("E12").Value = "  -8.50%  "
$ws.Range("E13").Value = "  -5.54%  "
$ws.Range("D14").Value = "2.777.73"
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("D15").Value = "56.449.74"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").Value = "'21.44"
$ws.Range("E16").Value = "  -4.62%  "
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").Value = "2.358.38"
$ws.Range("E18").Value = "  -5.91%  "
$ws.Range("D20").Value = "'4.05"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").Value = "'310.03"
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "'6.08"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D24").Value = "'65.09"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D27").Value = "2.463.41"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("D30").Value = "'174.89"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("E32").Value = "  -7.06%  "
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").Value = "'1.13"
$ws.Range("E34").Value = "  -6.82%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("E39").Value = "  -7.17%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").Value = "'35.47"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("E42").Value = "  -6.93%  "
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("D44").Value = "'4.83"
$ws.Range("E44").Value = "  -7.44%  "
$ws.Range("D45").Value = "'253.04"
$ws.Range("E45").Value = "  -9.87%  "
$ws.Range("D46").Value = "'0.567"
$ws.Range("E46").Value = "  -4.73%  "
$ws.Range("D47").Value = "'0.0905"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "'120.08"
$ws.Range("E48").Value = "  -9.06%  "
$ws.Range("D49").Value = "'0.0488"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("D51").Value = "'16.58"
$ws.Range("E51").Value = "  -6.93%  "
